$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 750
$ws.Range("I18").Value = 750
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 750
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -466

$ws.Range("H19").Value = 4304.3335
$ws.Range("I19").Value = 3849
$ws.Range("J19").Value = 4532
$ws.Range("K19").Value = 3849
$ws.Range("L19").Value = 4532
$ws.Range("M19").Value = -3674
$ws.Range("N19").Value = -4882

$ws.Range("H74").Value = 6156.6
$ws.Range("I74").Value = 6156.6
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 6156.6
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -5220.6

$ws.Range("H77").Value = 6156.6
$ws.Range("I77").Value = 6156.6
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 30783
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -26103

$ws.Range("H106").Value = 2406
$ws.Range("I106").Value = 2406
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 2406
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -1775

$ws.Range("H113").Value = 4924.75
$ws.Range("I113").Value = 2202.5
$ws.Range("J113").Value = 5832.1665
$ws.Range("K113").Value = 2202.5
$ws.Range("L113").Value = 5832.1665
$ws.Range("M113").Value = 1051.5
$ws.Range("N113").Value = -12340.1665

$ws.Range("H132").Value = 16785.305
$ws.Range("I132").Value = 18506.83
$ws.Range("J132").Value = 2668.8
$ws.Range("K132").Value = 55520.49000000001
$ws.Range("L132").Value = 8006.400000000001
$ws.Range("M132").Value = -52990.49000000001

$ws.Range("H137").Value = 11763.839
$ws.Range("I137").Value = 14324.75
$ws.Range("J137").Value = 2983.5715
$ws.Range("K137").Value = 42974.25
$ws.Range("L137").Value = 8950.7145
$ws.Range("M137").Value = -40424.25

$ws.Range("H138").Value = 27741.725
$ws.Range("I138").Value = 2267.3
$ws.Range("J138").Value = 53216.15
$ws.Range("K138").Value = 6801.900000000001
$ws.Range("L138").Value = 159648.45
$ws.Range("M138").Value = -1661.900000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24423.113
$ws.Range("I32").Value = 24423.113
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 24423.113
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -24136.113

$ws.Range("H63").Value = 4999.3335
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 4999.3335
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 4999.3335
$ws.Range("N63").Value = -6371.3335
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 4999.3335
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 4999.3335
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 24996.6675
$ws.Range("N66").Value = -31860.6675
$ws.Range("M66").ClearContents()

$ws.Range("H74").Value = 235345.97
$ws.Range("I74").Value = 273327.9
$ws.Range("J74").Value = 26445.25
$ws.Range("K74").Value = 273327.9
$ws.Range("L74").Value = 26445.25
$ws.Range("M74").Value = -272453.9
$ws.Range("N74").Value = -28193.25

$ws.Range("H77").Value = 235345.97
$ws.Range("I77").Value = 273327.9
$ws.Range("J77").Value = 26445.25
$ws.Range("K77").Value = 1366639.5
$ws.Range("L77").Value = 132226.25
$ws.Range("M77").Value = -1362271.5
$ws.Range("N77").Value = -140962.25

$ws.Range("H88").Value = 31000
$ws.Range("I88").Value = 12000
$ws.Range("J88").Value = 50000
$ws.Range("K88").Value = 12000
$ws.Range("L88").Value = 50000
$ws.Range("M88").Value = -11594
$ws.Range("N88").Value = -50812

$ws.Range("H91").Value = 31000
$ws.Range("I91").Value = 12000
$ws.Range("J91").Value = 50000
$ws.Range("K91").Value = 12000
$ws.Range("L91").Value = 50000
$ws.Range("M91").Value = -10596
$ws.Range("N91").Value = -52808

$ws.Range("H122").Value = 1803.05
$ws.Range("I122").Value = 1678.2354
$ws.Range("J122").Value = 2510.3333
$ws.Range("K122").Value = 5034.706200000001
$ws.Range("L122").Value = 7530.999899999999
$ws.Range("M122").Value = -2584.706200000001

$ws.Range("H132").Value = 1624.1476
$ws.Range("I132").Value = 1223.8572
$ws.Range("J132").Value = 2163
$ws.Range("K132").Value = 3671.5716
$ws.Range("L132").Value = 6489
$ws.Range("M132").Value = -1141.5716

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 9949.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 9949.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 9949.5
$ws.Range("N86").Value = -12195.5

$ws.Range("H89").Value = 9949.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 9949.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 49747.5
$ws.Range("N89").Value = -60979.5

$ws.Range("H99").Value = 2098.1538
$ws.Range("I99").Value = 1999.4445
$ws.Range("J99").Value = 2320.25
$ws.Range("K99").Value = 1999.4445
$ws.Range("L99").Value = 2320.25
$ws.Range("M99").Value = -501.4445000000001
$ws.Range("N99").Value = -5316.25

$ws.Range("H105").Value = 3678
$ws.Range("I105").Value = 1598
$ws.Range("J105").Value = 4926
$ws.Range("K105").Value = 1598
$ws.Range("L105").Value = 4926
$ws.Range("M105").Value = 149
$ws.Range("N105").Value = -8420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4547783
$ws.Range("I31").Value = 8334332
$ws.Range("J31").Value = 3923.9
$ws.Range("K31").Value = 8334332
$ws.Range("L31").Value = 3923.9
$ws.Range("M31").Value = -8334037
$ws.Range("N31").Value = -4513.9

$ws.Range("H34").Value = 4547783
$ws.Range("I34").Value = 8334332
$ws.Range("J34").Value = 3923.9
$ws.Range("K34").Value = 8334332
$ws.Range("L34").Value = 3923.9
$ws.Range("M34").Value = -8334130
$ws.Range("N34").Value = -4327.9

$ws.Range("H62").Value = 8537.462
$ws.Range("I62").Value = 8180.636
$ws.Range("J62").Value = 10500
$ws.Range("K62").Value = 8180.636
$ws.Range("L62").Value = 10500
$ws.Range("M62").Value = -7556.636

$ws.Range("H65").Value = 8537.462
$ws.Range("I65").Value = 8180.636
$ws.Range("J65").Value = 10500
$ws.Range("K65").Value = 40903.18
$ws.Range("L65").Value = 52500
$ws.Range("M65").Value = -37783.18

$ws.Range("H122").Value = 3579.8333
$ws.Range("I122").Value = 3615.8
$ws.Range("J122").Value = 3400
$ws.Range("K122").Value = 10847.4
$ws.Range("L122").Value = 10200
$ws.Range("M122").Value = -8397.400000000001
$ws.Range("N122").Value = -15100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 156.55556
$ws.Range("I26").Value = 152.85715
$ws.Range("J26").Value = 169.5
$ws.Range("K26").Value = 458.57145
$ws.Range("L26").Value = 508.5
$ws.Range("M26").Value = -170.57145
$ws.Range("N26").Value = -1084.5

$ws.Range("H68").Value = 4382.7144
$ws.Range("I68").Value = 1349
$ws.Range("J68").Value = 4566.5757
$ws.Range("K68").Value = 4047
$ws.Range("L68").Value = 13699.7271
$ws.Range("M68").Value = -3236
$ws.Range("N68").Value = -15321.7271

$ws.Range("H71").Value = 4382.7144
$ws.Range("I71").Value = 1349
$ws.Range("J71").Value = 4566.5757
$ws.Range("K71").Value = 12141
$ws.Range("L71").Value = 41099.1813
$ws.Range("M71").Value = -8085
$ws.Range("N71").Value = -49211.1813

$ws.Range("H107").Value = 978.6667
$ws.Range("I107").Value = 402.5
$ws.Range("J107").Value = 1143.2858
$ws.Range("K107").Value = 1207.5
$ws.Range("L107").Value = 3429.8574
$ws.Range("M107").Value = 712.5
$ws.Range("N107").Value = -7269.857400000001

$ws.Range("H129").Value = 2841.6
$ws.Range("I129").Value = 2086.4443
$ws.Range("J129").Value = 3974.3333
$ws.Range("K129").Value = 6259.3329
$ws.Range("L129").Value = 11922.9999
$ws.Range("M129").Value = -1259.3329

$ws.Range("H131").Value = 4388.727
$ws.Range("I131").Value = 6413.5
$ws.Range("J131").Value = 1959
$ws.Range("K131").Value = 19240.5
$ws.Range("L131").Value = 5877
$ws.Range("M131").Value = -14200.5
$ws.Range("N131").Value = -15957

$ws.Range("H138").Value = 2369.0833
$ws.Range("I138").Value = 1139.6
$ws.Range("J138").Value = 8516.5
$ws.Range("K138").Value = 3418.8
$ws.Range("L138").Value = 25549.5
$ws.Range("M138").Value = 1721.2
$ws.Range("N138").Value = -35829.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12996.75
$ws.Range("I70").Value = 11996.5
$ws.Range("J70").Value = 13997
$ws.Range("K70").Value = 11996.5
$ws.Range("L70").Value = 13997
$ws.Range("M70").Value = -11726.5
$ws.Range("N70").Value = -14537

$ws.Range("H73").Value = 12996.75
$ws.Range("I73").Value = 11996.5
$ws.Range("J73").Value = 13997
$ws.Range("K73").Value = 11996.5
$ws.Range("L73").Value = 13997
$ws.Range("M73").Value = -11060.5
$ws.Range("N73").Value = -15869

$ws.Range("H80").Value = 8098.5
$ws.Range("I80").Value = 1997.8
$ws.Range("J80").Value = 14199.2
$ws.Range("K80").Value = 1997.8
$ws.Range("L80").Value = 14199.2
$ws.Range("M80").Value = -999.8
$ws.Range("N80").Value = -16195.2

$ws.Range("H83").Value = 8098.5
$ws.Range("I83").Value = 1997.8
$ws.Range("J83").Value = 14199.2
$ws.Range("K83").Value = 9989
$ws.Range("L83").Value = 70996
$ws.Range("M83").Value = -4997
$ws.Range("N83").Value = -80980

$ws.Range("H97").Value = 673.25
$ws.Range("I97").Value = 747
$ws.Range("J97").Value = 452
$ws.Range("K97").Value = 747
$ws.Range("L97").Value = 452
$ws.Range("M97").Value = -251
$ws.Range("N97").Value = -1444

$ws.Range("H109").Value = 71095
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 71095
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 71095
$ws.Range("N109").Value = -73175

$ws.Range("H113").Value = 1900
$ws.Range("I113").Value = 1900
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1900
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 270
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 3229.4583
$ws.Range("I122").Value = 2981.6667
$ws.Range("J122").Value = 3972.8333
$ws.Range("K122").Value = 8945.000100000001
$ws.Range("L122").Value = 11918.4999
$ws.Range("M122").Value = -6495.000100000001
$ws.Range("N122").Value = -16818.4999

$ws.Range("H126").Value = 3519.182
$ws.Range("I126").Value = 2698.2
$ws.Range("J126").Value = 5278.4287
$ws.Range("K126").Value = 8094.599999999999
$ws.Range("L126").Value = 15835.2861
$ws.Range("M126").Value = -5624.599999999999
$ws.Range("N126").Value = -20775.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4702.9414
$ws.Range("I122").Value = 3425.9167
$ws.Range("J122").Value = 7767.8
$ws.Range("K122").Value = 10277.7501
$ws.Range("L122").Value = 23303.4
$ws.Range("M122").Value = -7827.750100000001

$ws.Range("H131").Value = 48930.25
$ws.Range("I131").Value = 50296
$ws.Range("J131").Value = 48475
$ws.Range("K131").Value = 50296
$ws.Range("L131").Value = 48475
$ws.Range("M131").Value = -45256
$ws.Range("N131").Value = -58555
